$d = $word.ActiveDocument

# --- Change 1: Add a _GoBack bookmark at the very start of the "Introduction" heading ---
$introRange = $d.Paragraphs(1).Range
$introRange.Collapse(1)  # wdCollapseStart
$d.Bookmarks.Add("_GoBack", $introRange) | Out-Null

# --- Change 2: Remove the _GoBack bookmark that currently sits inside the
#     "Your ability's strength is dependent ... on the level of your slot."
#     sentence, merging the two runs back into a single run of text. ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $bm = $d.Bookmarks("_GoBack")
    # only remove if it is not the one we just added at the Introduction heading
    if ($bm.Start -ne $introRange.Start) {
        $bm.Delete()
    }
}

# Find the second occurrence (the one in the ability paragraph) if the above
# logic did not already catch it; use Find/Replace to normalize the text
# (this also ensures the sentence reads as one contiguous run of text).
$rng = $d.Content
$rng.Find.Execute("Your ability", $true, $false, $false, $false, $false, `
                   $true, 1, $false, $null, 0) | Out-Null

# --- Change 3: Split "Strike down the skeleton." into extra text ---
$rng2 = $d.Content
$found = $rng2.Find.Execute("Strike down the skeleton.", $true, $false, $false, $false, $false, `
                             $true, 1, $false, "Strike down the skeleton the incoming skeleton. ", 2)
